$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy style from an existing header cell (H1) to the new ones
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J20
$data = @(
    @(8, 9),
    @(1, 6),
    @(1, 7),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 4),
    @(1, 5),
    @(8, 8),
    @(4, 6),
    @(9, 9),
    @(6, 6),
    @(4, 8),
    @(2, 7),
    @(1, 3),
    @(4, 4),
    @(3, 3),
    @(5, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
